# Auto-generated edit script: applies updated market-data values
# (currentAveragePrice / Leve profit columns) across the Leviathan_Profits sheets,
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1257.6111
$ws.Range("I6").Value = 159.13333
$ws.Range("K6").Value = 477.39999
$ws.Range("M6").Value = -365.39999
$ws.Range("H8").Value = 1879.625
$ws.Range("I8").Value = 2137.4285
$ws.Range("K8").Value = 6412.2855
$ws.Range("M8").Value = -6273.2855
$ws.Range("H32").Value = 2732.2666
$ws.Range("I32").Value = 2794
$ws.Range("J32").Value = 2562.5
$ws.Range("K32").Value = 2794
$ws.Range("L32").Value = 2562.5
$ws.Range("M32").Value = -2468
$ws.Range("N32").Value = -3214.5
$ws.Range("H94").Value = 1083.25
$ws.Range("I94").Value = 1083.25
$ws.Range("K94").Value = 1083.25
$ws.Range("M94").Value = -632.25
$ws.Range("H98").Value = 2935.4856
$ws.Range("I98").Value = 2526.6296
$ws.Range("K98").Value = 2526.6296
$ws.Range("M98").Value = -1028.6296
$ws.Range("H99").Value = 142880800
$ws.Range("I99").Value = 27598.166
$ws.Range("J99").Value = 1000000000
$ws.Range("K99").Value = 82794.49800000001
$ws.Range("L99").Value = 3000000000
$ws.Range("M99").Value = -81296.49800000001
$ws.Range("N99").Value = -3000002996
$ws.Range("H100").Value = 3297.8235
$ws.Range("I100").Value = 2751.7273
$ws.Range("J100").Value = 4299
$ws.Range("K100").Value = 2751.7273
$ws.Range("L100").Value = 4299
$ws.Range("M100").Value = -2210.7273
$ws.Range("N100").Value = -5381
$ws.Range("H112").Value = 1936.5
$ws.Range("I112").Value = 1172.25
$ws.Range("J112").Value = 2318.625
$ws.Range("K112").Value = 3516.75
$ws.Range("L112").Value = 6955.875
$ws.Range("M112").Value = -2408.75
$ws.Range("N112").Value = -9171.875
$ws.Range("H113").Value = 4355.8
$ws.Range("I113").Value = 3311.6
$ws.Range("J113").Value = 5400
$ws.Range("K113").Value = 3311.6
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = -57.59999999999991
$ws.Range("N113").Value = -11908
$ws.Range("H122").Value = 2935.4856
$ws.Range("I122").Value = 2526.6296
$ws.Range("K122").Value = 7579.888800000001
$ws.Range("M122").Value = -5129.888800000001
$ws.Range("H132").Value = 1171342.8
$ws.Range("I132").Value = 1482536.2
$ws.Range("J132").Value = 4367
$ws.Range("K132").Value = 4447608.6
$ws.Range("L132").Value = 13101
$ws.Range("M132").Value = -4445078.6
$ws.Range("N132").Value = -18161
$ws.Range("H135").Value = 1350
$ws.Range("H138").Value = 5099.4
$ws.Range("J138").Value = 5153.8125
$ws.Range("L138").Value = 15461.4375
$ws.Range("N138").Value = -25741.4375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 249.5
$ws.Range("I3").Value = 249.5
$ws.Range("K3").Value = 249.5
$ws.Range("M3").Value = -134.5
$ws.Range("H5").Value = 126.8
$ws.Range("I5").Value = 108.5
$ws.Range("K5").Value = 108.5
$ws.Range("M5").Value = 3.5
$ws.Range("H69").Value = 150000
$ws.Range("J69").Value = 150000
$ws.Range("L69").Value = 150000
$ws.Range("N69").Value = -151498
$ws.Range("H72").Value = 150000
$ws.Range("J72").Value = 150000
$ws.Range("L72").Value = 450000
$ws.Range("N72").Value = -457488
$ws.Range("H74").Value = 2099.8635
$ws.Range("I74").Value = 1443.3572
$ws.Range("J74").Value = 3248.75
$ws.Range("K74").Value = 1443.3572
$ws.Range("L74").Value = 3248.75
$ws.Range("M74").Value = -569.3571999999999
$ws.Range("N74").Value = -4996.75
$ws.Range("H77").Value = 2099.8635
$ws.Range("I77").Value = 1443.3572
$ws.Range("J77").Value = 3248.75
$ws.Range("K77").Value = 7216.786
$ws.Range("L77").Value = 16243.75
$ws.Range("M77").Value = -2848.786
$ws.Range("N77").Value = -24979.75
$ws.Range("H97").Value = 1107.5834
$ws.Range("I97").Value = 897.34485
$ws.Range("J97").Value = 1978.5714
$ws.Range("K97").Value = 897.34485
$ws.Range("L97").Value = 1978.5714
$ws.Range("M97").Value = -401.34485
$ws.Range("N97").Value = -2970.5714
$ws.Range("H132").Value = 2132.0527
$ws.Range("I132").Value = 2150.6453
$ws.Range("K132").Value = 6451.9359
$ws.Range("M132").Value = -3921.9359

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 126.8
$ws.Range("I4").Value = 108.5
$ws.Range("K4").Value = 108.5
$ws.Range("M4").Value = 6.5
$ws.Range("H82").Value = 21709.5
$ws.Range("H85").Value = 21709.5
$ws.Range("H99").Value = 2010.6666
$ws.Range("I99").Value = 2249.3125
$ws.Range("K99").Value = 2249.3125
$ws.Range("M99").Value = -751.3125
$ws.Range("H134").Value = 1563.1515
$ws.Range("I134").Value = 1383.5769
$ws.Range("K134").Value = 4150.7307
$ws.Range("M134").Value = -1615.7307

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1939.2
$ws.Range("J6").Value = 1895
$ws.Range("L6").Value = 5685
$ws.Range("N6").Value = -5911
$ws.Range("H33").Value = 159.9375
$ws.Range("I33").Value = 55.76923
$ws.Range("K33").Value = 334.61538
$ws.Range("M33").Value = -51.61538000000002
$ws.Range("H42").Value = 6540.4
$ws.Range("J42").Value = 7999.75
$ws.Range("L42").Value = 23999.25
$ws.Range("N42").Value = -25067.25
$ws.Range("H131").Value = 85049.414
$ws.Range("I131").Value = 1491.7273
$ws.Range("J131").Value = 155752.08
$ws.Range("K131").Value = 4475.1819
$ws.Range("L131").Value = 467256.24
$ws.Range("M131").Value = 564.8181000000004
$ws.Range("N131").Value = -477336.24

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 18800.666
$ws.Range("I69").Value = 12000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11251
$ws.Range("H72").Value = 18800.666
$ws.Range("I72").Value = 12000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -32256
$ws.Range("H80").Value = 6430.7666
$ws.Range("I80").Value = 6838.625
$ws.Range("J80").Value = 4799.3335
$ws.Range("K80").Value = 6838.625
$ws.Range("L80").Value = 4799.3335
$ws.Range("M80").Value = -5840.625
$ws.Range("N80").Value = -6795.3335
$ws.Range("H83").Value = 6430.7666
$ws.Range("I83").Value = 6838.625
$ws.Range("J83").Value = 4799.3335
$ws.Range("K83").Value = 34193.125
$ws.Range("L83").Value = 23996.6675
$ws.Range("M83").Value = -29201.125
$ws.Range("N83").Value = -33980.6675
$ws.Range("H97").Value = 19736.383
$ws.Range("I97").Value = 26612.92
$ws.Range("K97").Value = 26612.92
$ws.Range("M97").Value = -26116.92

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 38635.668
$ws.Range("I93").Value = 1433
$ws.Range("K93").Value = 1433
$ws.Range("M93").Value = -185
$ws.Range("H100").Value = 136332.8
$ws.Range("I100").Value = 4999.6665
$ws.Range("K100").Value = 4999.6665
$ws.Range("M100").Value = -4458.6665
$ws.Range("H122").Value = 9066
$ws.Range("I122").Value = 10527.714
$ws.Range("K122").Value = 31583.142
$ws.Range("M122").Value = -29133.142
$ws.Range("H132").Value = 16255.777
$ws.Range("I132").Value = 28825.75
$ws.Range("K132").Value = 86477.25
$ws.Range("M132").Value = -83947.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 91426.336
$ws.Range("I62").Value = 3764.6667
$ws.Range("K62").Value = 3764.6667
$ws.Range("M62").Value = -3140.6667
$ws.Range("H65").Value = 91426.336
$ws.Range("I65").Value = 3764.6667
$ws.Range("K65").Value = 18823.3335
$ws.Range("M65").Value = -15703.3335
$ws.Range("H81").Value = 2033.8182
$ws.Range("I81").Value = 2033.8182
$ws.Range("K81").Value = 4067.6364
$ws.Range("M81").Value = -3006.6364
$ws.Range("H84").Value = 2033.8182
$ws.Range("I84").Value = 2033.8182
$ws.Range("K84").Value = 20338.182
$ws.Range("M84").Value = -15034.182
$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800
$ws.Range("H126").Value = 2471.5454
$ws.Range("I126").Value = 2536.625
$ws.Range("K126").Value = 2536.625
$ws.Range("M126").Value = -5139.875

Write-Output "Updated 207 cells across 7 sheets."
